$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two data rows (4 and 5) have swapped their per-observation values
# (Id, Ost/Easting, Nord/Northing, Starttid, Sluttid) while the rest of the
# row (species, location name, dates, etc.) stayed identical between them.

# Row 4 <- values that used to be in row 5
$ws.Range("A4").Value = 131235752
$ws.Range("Q4").Value = 504836
$ws.Range("R4").Value = 6699938
$ws.Range("Z4").Value = "12:00"
$ws.Range("AB4").Value = "12:00"

# Row 5 <- values that used to be in row 4
$ws.Range("A5").Value = 131236495
$ws.Range("Q5").Value = 504860
$ws.Range("R5").Value = 6700261
$ws.Range("Z5").Value = "12:52"
$ws.Range("AB5").Value = "12:52"
